$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10
$ws.Range("E8").Value = "GIT UPDATE"

# Move/select the active cell to E8 to match the saved selection
$ws.Activate()
$ws.Range("E8").Select()
